{"js": "const replacements = [\n  [\"2025-08-18 Monday\", \"2025-08-19 Tuesday\"],\n  [\"662\u00f72=\", \"981\u00f73=\"],\n  [\"203\u00f74=\", \"132\u00f76=\"],\n  [\"660\u00f76=\", \"702\u00f75=\"],\n  [\"231\u00f79=\", \"465\u00f79=\"],\n  [\"158\u00f77=\", \"792\u00f74=\"],\n  [\"664\u00f78=\", \"774\u00f77=\"],\n  [\"658\u00f76=\", \"201\u00f74=\"],\n  [\"606\u00f78=\", \"588\u00f74=\"],\n  [\"727\u00f78=\", \"790\u00f77=\"],\n  [\"534\u00f76=\", \"129\u00f74=\"],\n  [\"671\u00f73=\", \"216\u00f74=\"],\n  [\"554\u00f77=\", \"187\u00f78=\"],\n  [\"396\u00f73=\", \"168\u00f73=\"],\n  [\"432\u00f74=\", \"697\u00f76=\"],\n  [\"911\u00f79=\", \"917\u00f79=\"],\n  [\"518\u00f76=\", \"177\u00f76=\"],\n  [\"631\u00f76=\", \"643\u00f76=\"],\n  [\"195\u00f74=\", \"483\u00f78=\"],\n  [\"884\u00f74=\", \"764\u00f72=\"],\n  [\"678\u00f74=\", \"800\u00f75=\"],\n  [\"395\u00f73=\", \"221\u00f77=\"],\n  [\"970\u00f79=\", \"971\u00f78=\"],\n  [\"834\u00f76=\", \"694\u00f75=\"],\n  [\"793\u00f72=\", \"507\u00f77=\"],\n  [\"845\u00f78=\", \"717\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-18 Monday\", \"2025-08-19 Tuesday\"),\n    @(\"662\u00f72=\", \"981\u00f73=\"),\n    @(\"203\u00f74=\", \"132\u00f76=\"),\n    @(\"660\u00f76=\", \"702\u00f75=\"),\n    @(\"231\u00f79=\", \"465\u00f79=\"),\n    @(\"158\u00f77=\", \"792\u00f74=\"),\n    @(\"664\u00f78=\", \"774\u00f77=\"),\n    @(\"658\u00f76=\", \"201\u00f74=\"),\n    @(\"606\u00f78=\", \"588\u00f74=\"),\n    @(\"727\u00f78=\", \"790\u00f77=\"),\n    @(\"534\u00f76=\", \"129\u00f74=\"),\n    @(\"671\u00f73=\", \"216\u00f74=\"),\n    @(\"554\u00f77=\", \"187\u00f78=\"),\n    @(\"396\u00f73=\", \"168\u00f73=\"),\n    @(\"432\u00f74=\", \"697\u00f76=\"),\n    @(\"911\u00f79=\", \"917\u00f79=\"),\n    @(\"518\u00f76=\", \"177\u00f76=\"),\n    @(\"631\u00f76=\", \"643\u00f76=\"),\n    @(\"195\u00f74=\", \"483\u00f78=\"),\n    @(\"884\u00f74=\", \"764\u00f72=\"),\n    @(\"678\u00f74=\", \"800\u00f75=\"),\n    @(\"395\u00f73=\", \"221\u00f77=\"),\n    @(\"970\u00f79=\", \"971\u00f78=\"),\n    @(\"834\u00f76=\", \"694\u00f75=\"),\n    @(\"793\u00f72=\", \"507\u00f77=\"),\n    @(\"845\u00f78=\", \"717\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
